# Update res_bus/vm_pu.xlsx values for the "case with 380 kV" power-flow re-run.
# The slack/reference bus setpoint (column B) moves from 1.05 pu to 1.02 pu,
# and all other bus voltage magnitudes (columns C-F, I-N) are refreshed with the
# newly computed results for rows 2-25 (buses 0-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (bus 0)
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.045457592443495
$ws.Range("D2").Value = 1.04557819089346
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.052934288542099
$ws.Range("I2").Value = 1.042587921345068
$ws.Range("J2").Value = 1.050517222926146
$ws.Range("K2").Value = 1.048346037851956
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.055681648964118
$ws.Range("N2").Value = 1.052009077878711

# Row 3 (bus 1)
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04659938035503
$ws.Range("D3").Value = 1.046431228425557
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.054280666365417
$ws.Range("I3").Value = 1.04292658005496
$ws.Range("J3").Value = 1.05130573138442
$ws.Range("K3").Value = 1.04901043596468
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.056839604052908
$ws.Range("N3").Value = 1.052798706109439

# Row 4 (bus 2)
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04733762000706
$ws.Range("D4").Value = 1.046982705143404
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.055151713641695
$ws.Range("I4").Value = 1.043144267227801
$ws.Range("J4").Value = 1.051814829565823
$ws.Range("K4").Value = 1.049439209346241
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.057588191283284
$ws.Range("N4").Value = 1.053308527268654

# Row 5 (bus 3)
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.047647840579523
$ws.Range("D5").Value = 1.047214428021014
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.055517869004017
$ws.Range("I5").Value = 1.043235436901588
$ws.Range("J5").Value = 1.052028587604522
$ws.Range("K5").Value = 1.049619194427273
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.057902735264406
$ws.Range("N5").Value = 1.053522588868286

# Row 6 (bus 4)
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047699920095646
$ws.Range("D6").Value = 1.047253328416713
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.05557934625203
$ws.Range("I6").Value = 1.043250724416772
$ws.Range("J6").Value = 1.052064462906514
$ws.Range("K6").Value = 1.049649398833162
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.057955539147885
$ws.Range("N6").Value = 1.053558515117322

# Row 7 (bus 5)
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.047341765719081
$ws.Range("D7").Value = 1.046985801898354
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.055156606353385
$ws.Range("I7").Value = 1.043145486800346
$ws.Range("J7").Value = 1.051817686858219
$ws.Range("K7").Value = 1.049441615380455
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.057592394869249
$ws.Range("N7").Value = 1.053311388618733

# Row 8 (bus 6)
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045843584909744
$ws.Range("D8").Value = 1.045866581972351
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.053389335816315
$ws.Range("I8").Value = 1.04270267272584
$ws.Range("J8").Value = 1.050783935565578
$ws.Range("K8").Value = 1.048570810170775
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.056073129038953
$ws.Range("N8").Value = 1.052276169280683

# Row 9 (bus 7)
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04319912057876
$ws.Range("D9").Value = 1.043890544989754
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.050273910587716
$ws.Range("I9").Value = 1.041911265107914
$ws.Range("J9").Value = 1.048953706272953
$ws.Range("K9").Value = 1.047027595832306
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.053390617464929
$ws.Range("N9").Value = 1.050443340852533

# Row 10 (bus 8)
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041433010307072
$ws.Range("D10").Value = 1.04257056436606
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.048195928988961
$ws.Range("I10").Value = 1.041376151760264
$ws.Range("J10").Value = 1.047727675983776
$ws.Range("K10").Value = 1.045992850186421
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.051598511480438
$ws.Range("N10").Value = 1.049215569459717

# Row 11 (bus 9)
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040667493515253
$ws.Range("D11").Value = 1.041998364436123
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.047295854974237
$ws.Range("I11").Value = 1.041142652292423
$ws.Range("J11").Value = 1.047195379453442
$ws.Range("K11").Value = 1.045543372572504
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.05082158191389
$ws.Range("N11").Value = 1.048682517007255

# Row 12 (bus 10)
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040383026866614
$ws.Range("D12").Value = 1.041785726700365
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.046961479625378
$ws.Range("I12").Value = 1.041055650177964
$ws.Range("J12").Value = 1.046997446250598
$ws.Range("K12").Value = 1.045376201101825
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.050532851977021
$ws.Range("N12").Value = 1.04848430271656

# Row 13 (bus 11)
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040444051382574
$ws.Range("D13").Value = 1.041831342633149
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.047033206492609
$ws.Range("I13").Value = 1.041074324658658
$ws.Range("J13").Value = 1.047039913330626
$ws.Range("K13").Value = 1.045412069720366
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.050594792081396
$ws.Range("N13").Value = 1.048526830104712

# Row 14 (bus 12)
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040643981868388
$ws.Range("D14").Value = 1.041980789723954
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.047268216399081
$ws.Range("I14").Value = 1.041135466183429
$ws.Range("J14").Value = 1.047179022619813
$ws.Range("K14").Value = 1.045529558529367
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.05079771835672
$ws.Range("N14").Value = 1.048666136945046

# Row 15 (bus 13)
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.040767149711487
$ws.Range("D15").Value = 1.042072856095041
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.047413007319824
$ws.Range("I15").Value = 1.041173101690818
$ws.Range("J15").Value = 1.047264703962309
$ws.Range("K15").Value = 1.045601918672775
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.050922728797322
$ws.Range("N15").Value = 1.048751939964876

# Row 16 (bus 14)
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041483798498554
$ws.Range("D16").Value = 1.042608525820516
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.048255657432349
$ws.Range("I16").Value = 1.04139161052957
$ws.Range("J16").Value = 1.047762972750742
$ws.Range("K16").Value = 1.046022650386867
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.051650053755197
$ws.Range("N16").Value = 1.049250916352141

# Row 17 (bus 15)
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041933122900092
$ws.Range("D17").Value = 1.04294436526508
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.048784148283302
$ws.Range("I17").Value = 1.041528195050305
$ws.Range("J17").Value = 1.048075142989318
$ws.Range("K17").Value = 1.046286181641093
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.052106032746453
$ws.Range("N17").Value = 1.049563529908262

# Row 18 (bus 16)
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.04219513119226
$ws.Range("D18").Value = 1.04314019309318
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.049092379902922
$ws.Range("I18").Value = 1.041607689656691
$ws.Range("J18").Value = 1.048257089966806
$ws.Range("K18").Value = 1.046439757466236
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.05237190749034
$ws.Range("N18").Value = 1.049745735271328

# Row 19 (bus 17)
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042284456599605
$ws.Range("D19").Value = 1.043206954900673
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.04919747425316
$ws.Range("I19").Value = 1.041634765972103
$ws.Range("J19").Value = 1.048319106057054
$ws.Range("K19").Value = 1.046492099569704
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.052462548721039
$ws.Range("N19").Value = 1.049807839431537

# Row 20 (bus 18)
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041884922422765
$ws.Range("D20").Value = 1.04290833925434
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.048727449153903
$ws.Range("I20").Value = 1.041513558706575
$ws.Range("J20").Value = 1.048041664222063
$ws.Range("K20").Value = 1.04625792146227
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.0520571198562
$ws.Range("N20").Value = 1.049530003597318

# Row 21 (bus 19)
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040585110646151
$ws.Range("D21").Value = 1.041936783972872
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.047199013186923
$ws.Range("I21").Value = 1.041117468988564
$ws.Range("J21").Value = 1.047138064337729
$ws.Range("K21").Value = 1.04549496696155
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.050737965618318
$ws.Range("N21").Value = 1.048625120497503

# Row 22 (bus 20)
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039767174003342
$ws.Range("D22").Value = 1.041325365653892
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.046237746589419
$ws.Range("I22").Value = 1.040866869032141
$ws.Range("J22").Value = 1.046568692970023
$ws.Range("K22").Value = 1.045014019910386
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.049907728026285
$ws.Range("N22").Value = 1.048054940557147

# Row 23 (bus 21)
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040200844267667
$ws.Range("D23").Value = 1.041649543808147
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.046747359771424
$ws.Range("I23").Value = 1.040999865201111
$ws.Range("J23").Value = 1.046870645804538
$ws.Range("K23").Value = 1.045269097686644
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.050347932677103
$ws.Range("N23").Value = 1.048357322199322

# Row 24 (bus 22)
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041906702374879
$ws.Range("D24").Value = 1.042924618047971
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.048753069136007
$ws.Range("I24").Value = 1.041520172774083
$ws.Range("J24").Value = 1.048056792257653
$ws.Range("K24").Value = 1.046270691445113
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.052079221763454
$ws.Range("N24").Value = 1.049545153116453

# Row 25 (bus 23)
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043883321340797
$ws.Range("D25").Value = 1.044401856576896
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.0510794913734
$ws.Range("I25").Value = 1.042117183623899
$ws.Range("J25").Value = 1.049427894692412
$ws.Range("K25").Value = 1.047427595807153
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.054084763301809
$ws.Range("N25").Value = 1.050918202673937
